{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = (p.text || \"\").trim();\n  if (p.style === \"Author\" && (text === \"Ben Jarman\" || text === \"Helen Fair\")) {\n    p.delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text.Trim()\n    $styleName = $para.Range.Style.NameLocal\n    if ($styleName -eq \"Author\" -and ($text -eq \"Ben Jarman\" -or $text -eq \"Helen Fair\")) {\n        $para.Range.Delete()\n    }\n}\n"}
